$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''54.326.29'
$ws.Range('E2').Value = '''  -6.79%  '
$ws.Range('D3').Value = '''2.437.63'
$ws.Range('E3').Value = '''  -9.42%  '
$ws.Range('E4').Value = '''  +0.05%  '
$ws.Range('D5').Value = '''467.45'
$ws.Range('E5').Value = '''  -6.25%  '
$ws.Range('D6').Value = '''130.64'
$ws.Range('E6').Value = '''  -5.56%  '
$ws.Range('E7').Value = '''  +0.22%  '
$ws.Range('D8').Value = '''0.492'
$ws.Range('E8').Value = '''  -6.45%  '
$ws.Range('D9').Value = '''2.433.71'
$ws.Range('E9').Value = '''  -9.86%  '
$ws.Range('D10').Value = '''0.0948'
$ws.Range('E10').Value = '''  -9.02%  '
$ws.Range('E11').Value = '''  -12.05%  '
$ws.Range('D12').Value = '''0.314'
$ws.Range('E12').Value = '''  -8.99%  '
$ws.Range('D14').Value = '''2.870.49'
$ws.Range('E14').Value = '''  -9.34%  '
$ws.Range('D15').Value = '''54.461.57'
$ws.Range('E15').Value = '''  -6.69%  '
$ws.Range('E16').Value = '''  -0.19%  '
$ws.Range('D17').Value = '''19.58'
$ws.Range('E17').Value = '''  -8.10%  '
$ws.Range('D18').Value = '''2.446.93'
$ws.Range('E18').Value = '''  -9.45%  '
$ws.Range('E19').Value = '''  -10.12%  '
$ws.Range('D20').Value = '''311.46'
$ws.Range('D21').Value = '''9.53'
$ws.Range('E21').Value = '''  -12.69%  '
$ws.Range('D22').Value = '''0.999'
$ws.Range('E22').Value = '''  +0.38%  '
$ws.Range('D23').Value = '''5.69'
$ws.Range('E23').Value = '''  +0.96%  '
$ws.Range('D24').Value = '''5.39'
$ws.Range('E24').Value = '''  -12.84%  '
$ws.Range('D25').Value = '''56.29'
$ws.Range('E25').Value = '''  -10.04%  '
$ws.Range('E26').Value = '''  +0.77%  '
$ws.Range('E27').Value = '''  -8.70%  '
$ws.Range('D28').Value = '''2.551.90'
$ws.Range('E28').Value = '''  -9.38%  '
$ws.Range('D29').Value = '''0.155'
$ws.Range('E29').Value = '''  -8.32%  '
$ws.Range('E30').Value = '''  -3.20%  '
$ws.Range('E31').Value = '''  +0.07%  '
$ws.Range('D32').Value = '''0.0₃0709'
$ws.Range('E32').Value = '''  -12.85%  '
$ws.Range('D33').Value = '''145.63'
$ws.Range('E33').Value = '''  -2.94%  '
$ws.Range('E34').Value = '''  -6.39%  '
$ws.Range('E35').Value = '''  -9.51%  '
$ws.Range('E36').Value = '''  -6.01%  '
$ws.Range('D37').Value = '''3.55'
$ws.Range('E37').Value = '''  -14.52%  '
$ws.Range('E38').Value = '''  -5.18%  '
$ws.Range('D39').Value = '''0.792'
$ws.Range('E39').Value = '''  -14.56%  '
$ws.Range('D41').Value = '''32.76'
$ws.Range('E41').Value = '''  -6.94%  '
$ws.Range('D42').Value = '''0.596'
$ws.Range('E42').Value = '''  +0.93%  '
$ws.Range('E43').Value = '''  -5.48%  '
$ws.Range('D44').Value = '''3.25'
$ws.Range('E44').Value = '''  -8.17%  '
$ws.Range('D45').Value = '''10.11'
$ws.Range('E45').Value = '''  -2.40%  '
$ws.Range('D46').Value = '''1.23'
$ws.Range('E46').Value = '''  -9.72%  '
$ws.Range('D47').Value = '''1.929.92'
$ws.Range('E47').Value = '''  -11.09%  '
$ws.Range('D48').Value = '''0.0883'
$ws.Range('E48').Value = '''  -0.05%  '
$ws.Range('D49').Value = '''0.0216'
$ws.Range('E49').Value = '''  -3.88%  '
$ws.Range('D50').Value = '''231.80'
$ws.Range('E50').Value = '''  +6.14%  '
$ws.Range('D51').Value = '''16.58'
$ws.Range('E51').Value = '''  -10.78%  '
